# Applies the cryptos list refresh described in the commit:
# "Updated cryptos list on Thu Dec 28 08:55:51 UTC 2023 with GitHub Actions"
#
# Column D holds price text that often looks numeric (e.g. "2.18", "327.75").
# The source workbook stores these as plain text (inlineStr), so when we assign
# such a value Excel would otherwise auto-convert it to a Double. Set-TextValue
# forces the cell to Text format just long enough to assign the literal string,
# then restores the cell's original style so no visible formatting changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $cell = $ws.Range($addr)
    $savedStyle = $cell.Style
    $cell.NumberFormat = '@'
    $cell.Value = $text
    $cell.Style = $savedStyle
}

$ws.Range('D2').Value = '43.124.11'
$ws.Range('D3').Value = '2.389.02'
$ws.Range('E3').Value = '  +6.19%  '
$ws.Range('E4').Value = '  -0.53%  '
Set-TextValue 'D5' '327.75'
$ws.Range('E5').Value = '  +10.81%  '
Set-TextValue 'D6' '105.68'
$ws.Range('E6').Value = '  -6.86%  '
Set-TextValue 'D7' '0.644'
$ws.Range('E7').Value = '  +2.44%  '
Set-TextValue 'D9' '0.657'
$ws.Range('E9').Value = '  +8.40%  '
Set-TextValue 'D10' '41.84'
$ws.Range('E10').Value = '  -5.11%  '
Set-TextValue 'D11' '0.0937'
$ws.Range('E11').Value = '  +1.50%  '
$ws.Range('E12').Value = '  -3.03%  '
$ws.Range('E13').Value = '  -1.39%  '
$ws.Range('E14').Value = '  +13.95%  '
$ws.Range('E15').Value = '  +2.09%  '
$ws.Range('D16').Value = '2.750.11'
$ws.Range('E16').Value = '  +6.32%  '
$ws.Range('D17').Value = '2.382.49'
$ws.Range('E17').Value = '  +6.07%  '
$ws.Range('D18').Value = '43.131.08'
$ws.Range('E18').Value = '  +0.86%  '
$ws.Range('E19').Value = '  +9.60%  '
$ws.Range('E20').Value = '  +2.33%  '
Set-TextValue 'D21' '76.74'
$ws.Range('E21').Value = '  +2.87%  '
Set-TextValue 'D22' '3.63'
$ws.Range('E22').Value = '  +3.39%  '
Set-TextValue 'D23' '275.18'
$ws.Range('E23').Value = '  +10.51%  '
Set-TextValue 'D25' '9.58'
$ws.Range('E25').Value = '  +6.85%  '
Set-TextValue 'D26' '11.78'
$ws.Range('E26').Value = '  +2.05%  '
$ws.Range('E27').Value = '  -0.01%  '
Set-TextValue 'D28' '23.13'
$ws.Range('E28').Value = '  +4.98%  '
Set-TextValue 'D29' '37.80'
$ws.Range('E29').Value = '  +0.20%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D30' '2.18'
$ws.Range('E30').Value = '  -1.45%  '
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D31' '174.88'
$ws.Range('E31').Value = '  -0.54%  '
$ws.Range('E32').Value = '  +1.33%  '
Set-TextValue 'D33' '0.0929'
$ws.Range('E33').Value = '  +4.69%  '
Set-TextValue 'D34' '5.88'
$ws.Range('E34').Value = '  +3.11%  '
$ws.Range('E35').Value = '  +5.26%  '
$ws.Range('E36').Value = '  -3.07%  '
$ws.Range('E37').Value = '  -0.62%  '
Set-TextValue 'D38' '0.0367'
$ws.Range('E38').Value = '  -2.31%  '
Set-TextValue 'D39' '0.107'
$ws.Range('E39').Value = '  +2.38%  '
Set-TextValue 'D40' '2.81'
$ws.Range('E40').Value = '  +16.14%  '
$ws.Range('E41').Value = '  +19.13%  '
Set-TextValue 'D42' '0.234'
$ws.Range('E42').Value = '  +1.19%  '
Set-TextValue 'D43' '69.74'
$ws.Range('E43').Value = '  -3.35%  '
Set-TextValue 'D44' '121.57'
$ws.Range('E44').Value = '  +15.41%  '
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('B46').Value = 'Celestia'
$ws.Range('C46').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue 'D46' '12.49'
$ws.Range('E46').Value = '  -0.07%  '
$ws.Range('B47').Value = 'BitcoinSV'
$ws.Range('C47').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue 'D47' '90.49'
$ws.Range('E47').Value = '  +54.38%  '
Set-TextValue 'D48' '9.39'
$ws.Range('E48').Value = '  +8.86%  '
$ws.Range('E49').Value = '  +0.51%  '
Set-TextValue 'D50' '1.31'
$ws.Range('E50').Value = '  +1.00%  '
$ws.Range('D51').Value = '1.592.42'
$ws.Range('E51').Value = '  +10.23%  '
